# Replicate the previous "Hate Crime" block-append step, this time for
# Offense = 'Hate Crime - Statutory Rape', with data only for
# Date = 'sum2014' and 'sum2015' (mirrors HateCrimes.xls source, which had
# no sum2013 column for this offense).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offense = "Hate Crime - Statutory Rape"

# Column A cycles through the 9 sector-of-institution values (in sheet order).
$sectors = @(
    "Public, 4-year or above",
    "Private nonprofit, 4-year or above",
    "Private for-profit, 4-year or above",
    "Public, 2-year",
    "Private nonprofit, 2-year",
    "Private for-profit, 2-year",
    "Public, less-than 2-year",
    "Private nonprofit, less-than 2-year",
    "Private for-profit, less-than 2-year"
)

# Column B cycles through the 4 reporting locations.
$locations = @(
    "On Campus (excluding Residence Halls)",
    "On Campus (Residence Halls)",
    "Non-Campus",
    "Public Property"
)

# Column D only has the two years that existed in the source data for this
# offense (no sum2013).
$dates = @("sum2014", "sum2015")

$startRow = 2189
$row = $startRow

foreach ($loc in $locations) {
    foreach ($date in $dates) {
        foreach ($sector in $sectors) {
            $ws.Rows.Item($row).Insert()

            $ws.Cells.Item($row, 1).Value = "'" + $sector
            $ws.Cells.Item($row, 2).Value = "'" + $loc
            $ws.Cells.Item($row, 3).Value = $offense
            $ws.Cells.Item($row, 4).Value = "'" + $date
            $ws.Cells.Item($row, 5).Value = 0

            $row = $row + 1
        }
    }
}

$endRow = $row - 1

$ws.Range("G2195").Select()
